$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 "Enterprises by employment size band" - update source link and dates
$ws.Range("B8").Value = "<a href='https://www.nomisweb.co.uk/datasets/idbrent'>ONS UK Business Counts</a>"
$ws.Range("C8").Value = "Mar 2023 (27/09/23)"
$ws.Range("D8").Value = "Mar 2024 (09/24)"

# Row 9 "Enterprises by employment industry" - update the same released/latest period dates
$ws.Range("C9").Value = "Mar 2023 (27/09/23)"
$ws.Range("D9").Value = "Mar 2024 (09/24)"

# Move active selection from B5 to B6, as reflected in the saved view state
$ws.Range("B6").Select()
